$d = $word.ActiveDocument

function FindSubRange($startPos, $endPos, $searchText) {
    $r = $d.Range($startPos, $endPos)
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { throw "text not found: $searchText" }
    return $r
}

# ---------------------------------------------------------------------------
# 1. Locate the fixed paragraphs that hold a single "value" run of text.
# ---------------------------------------------------------------------------
$pObjetivos = $d.Paragraphs.Item(6)
$pDocente   = $d.Paragraphs.Item(8)
$pResumido  = $d.Paragraphs.Item(10)
$pPrograma  = $d.Paragraphs.Item(12)
$pAvaliacao = $d.Paragraphs.Item(14)
$pBiblio    = $d.Paragraphs.Item(16)

# The "Avaliação" paragraph holds three separate value spans, interleaved
# with bold labels ("Método: ", "Critério: ", "Norma de recuperação: ").
$rMetodo  = FindSubRange $pAvaliacao.Range.Start $pAvaliacao.Range.End "Aulas expositivas; microcomputadores; seminários; visitas técnicas."
$rCriterio = FindSubRange $pAvaliacao.Range.Start $pAvaliacao.Range.End "Média ponderada de notas de provas e seminários."
$rNorma   = FindSubRange $pAvaliacao.Range.Start $pAvaliacao.Range.End "Prova única com nota igual ou superior a 5,0 (cinco)."

# ---------------------------------------------------------------------------
# 2. Capture the *current* (pre-edit) text of every slot, manual line breaks
#    included (Word represents a <w:br/> inside Range.Text as Chr(11)).
# ---------------------------------------------------------------------------
$txtObjetivos = $pObjetivos.Range.Text
$txtDocente   = $pDocente.Range.Text
$txtResumido  = $pResumido.Range.Text
$txtPrograma  = $pPrograma.Range.Text
$txtMetodo    = $rMetodo.Text
$txtCriterio  = $rCriterio.Text
$txtNorma     = $rNorma.Text
$txtBiblio    = $pBiblio.Range.Text

# Paragraph Range.Text includes the trailing paragraph mark (Chr(13)); strip
# it so we can reassign the body text without destroying the paragraph mark.
function StripParaMark($s) {
    if ($s.Length -gt 0 -and [int][char]$s[$s.Length - 1] -eq 13) {
        return $s.Substring(0, $s.Length - 1)
    }
    return $s
}

$txtObjetivos = StripParaMark $txtObjetivos
$txtDocente   = StripParaMark $txtDocente
$txtResumido  = StripParaMark $txtResumido
$txtPrograma  = StripParaMark $txtPrograma
$txtBiblio    = StripParaMark $txtBiblio

# ---------------------------------------------------------------------------
# 3. Re-assign the captured texts to their new (rotated) homes.
#    new(slot) = old(next slot in the cycle)
#    Objetivos <- Resumido <- Programa <- Metodo <- Criterio <- Norma <- Biblio <- Docente <- Objetivos
#
#    $rMetodo / $rCriterio / $rNorma are plain Range objects (not live
#    Paragraph refs), so they do NOT auto-shift when text to their *left*
#    changes length. Writing them back in right-to-left document order
#    keeps every not-yet-written Range's Start/End valid at the moment it
#    is used (an edit only ever shifts what comes after it).
# ---------------------------------------------------------------------------
$rNorma.Text           = $txtBiblio
$rCriterio.Text        = $txtNorma
$rMetodo.Text          = $txtCriterio

# These are live Paragraph objects: their .Range re-resolves on access, so
# order does not matter here.
$pObjetivos.Range.Text = $txtResumido
$pDocente.Range.Text   = $txtObjetivos
$pResumido.Range.Text  = $txtPrograma
$pPrograma.Range.Text  = $txtMetodo
$pBiblio.Range.Text    = $txtDocente

Write-Host "done"
